# Update "Horarios Línea 141" workbook with the latest scrape results.
# New scrape timestamp:
$newTime = "02:13:28"

$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912": add 2 new schedule rows, bump header counters ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 7"

$ws1.Cells.Item(11, 1).Value = $newTime
$ws1.Cells.Item(11, 2).Value = "03:55"
$ws1.Cells.Item(11, 3).Value = "14_ABASTO"
$ws1.Cells.Item(11, 4).Value = 102
$ws1.Cells.Item(11, 5).Value = "LP1912"

$ws1.Cells.Item(12, 1).Value = $newTime
$ws1.Cells.Item(12, 2).Value = "04:01"
$ws1.Cells.Item(12, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(12, 4).Value = 108
$ws1.Cells.Item(12, 5).Value = "LP1912"

# --- Sheet "LP1912-215": just bump the "last update" timestamp ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

# --- Sheet "6203-6173": just bump the "last update" timestamp ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
